# Update crypto price/volume data per the Sat Jul  6 05:00:09 UTC 2024 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet stores every Price/Volume cell as text (e.g. "500.88", "  +5.86%  ").
# Plain decimal numbers would otherwise auto-convert to the Number type on assignment,
# so mark those specific Price cells as Text first to keep them stored like their neighbours.
$textFormatCells = @("D5", "D6", "D8", "D9", "D11", "D14", "D16", "D18", "D19", "D20", "D21", "D22", "D24", "D25", "D28", "D29", "D32", "D33", "D34", "D36", "D38", "D41", "D42", "D45", "D47", "D48", "D50", "D51")
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '56.072.36'
$ws.Range("E2").Value = '  +3.59%  '

$ws.Range("D3").Value = '2.963.55'
$ws.Range("E3").Value = '  +2.95%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = '500.19'
$ws.Range("E5").Value = '  +6.14%  '

$ws.Range("D6").Value = '133.87'
$ws.Range("E6").Value = '  +6.46%  '

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("D8").Value = '0.427'
$ws.Range("E8").Value = '  +5.99%  '

$ws.Range("D9").Value = '7.32'
$ws.Range("E9").Value = '  +10.48%  '

$ws.Range("E10").Value = '  +9.44%  '

$ws.Range("D11").Value = '0.349'
$ws.Range("E11").Value = '  +4.91%  '

$ws.Range("E12").Value = '  +3.22%  '

$ws.Range("D13").Value = '3.472.28'
$ws.Range("E13").Value = '  +2.79%  '

$ws.Range("D14").Value = '25.18'
$ws.Range("E14").Value = '  +10.17%  '

$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '56.135.75'
$ws.Range("E15").Value = '  +3.62%  '

$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '0.0000149'
$ws.Range("E16").Value = '  +11.53%  '

$ws.Range("D17").Value = '2.967.02'
$ws.Range("E17").Value = '  +2.67%  '

$ws.Range("D18").Value = '5.66'
$ws.Range("E18").Value = '  +9.22%  '

$ws.Range("D19").Value = '12.26'
$ws.Range("E19").Value = '  +6.32%  '

$ws.Range("D20").Value = '7.70'
$ws.Range("E20").Value = '  +8.24%  '

$ws.Range("D21").Value = '320.60'
$ws.Range("E21").Value = '  +4.66%  '

$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  +0.10%  '

$ws.Range("E23").Value = '  +4.53%  '

$ws.Range("D24").Value = '61.66'
$ws.Range("E24").Value = '  +3.99%  '

$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.12%  '

$ws.Range("E26").Value = '  +5.56%  '

$ws.Range("D27").Value = '0.0₃0878'
$ws.Range("E27").Value = '  +7.23%  '

$ws.Range("D28").Value = '6.42'
$ws.Range("E28").Value = '  +4.46%  '

$ws.Range("D29").Value = '6.71'
$ws.Range("E29").Value = '  +9.23%  '

$ws.Range("E30").Value = '  +2.87%  '

$ws.Range("E31").Value = '  +7.84%  '

$ws.Range("D32").Value = '20.33'
$ws.Range("E32").Value = '  +6.82%  '

$ws.Range("D33").Value = '157.99'
$ws.Range("E33").Value = '  +13.26%  '

$ws.Range("D34").Value = '4.42'
$ws.Range("E34").Value = '  +4.43%  '

$ws.Range("E35").Value = '  +3.28%  '

$ws.Range("D36").Value = '5.51'
$ws.Range("E36").Value = '  +1.26%  '

$ws.Range("E37").Value = '  +8.16%  '

$ws.Range("D38").Value = '22.81'
$ws.Range("E38").Value = '  -0.32%  '

$ws.Range("D39").Value = '2.996.98'
$ws.Range("E39").Value = '  +2.97%  '

$ws.Range("E40").Value = '  -0.05%  '

$ws.Range("D41").Value = '36.14'
$ws.Range("E41").Value = '  +4.40%  '

$ws.Range("D42").Value = '0.636'
$ws.Range("E42").Value = '  +6.92%  '

$ws.Range("D43").Value = '2.235.44'
$ws.Range("E43").Value = '  +8.77%  '

$ws.Range("E44").Value = '  +5.49%  '

$ws.Range("D45").Value = '0.966'
$ws.Range("E45").Value = '  +0.17%  '

$ws.Range("E46").Value = '  +2.77%  '

$ws.Range("D47").Value = '1.90'
$ws.Range("E47").Value = '  +18.91%  '

$ws.Range("D48").Value = '5.72'
$ws.Range("E48").Value = '  +7.84%  '

$ws.Range("E49").Value = '  +9.53%  '

$ws.Range("D50").Value = '18.81'
$ws.Range("E50").Value = '  +4.24%  '

$ws.Range("D51").Value = '0.0862'
$ws.Range("E51").Value = '  +7.91%  '
